# Apply translation/localization updates + new dialogue rows for old_wizard
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New labelled cells in column H ("main menu" / "strings") ---
$ws.Range("H2").Value = "main menu"
$ws.Range("H3").Value = "strings"

# --- New dialogue rows (6-9): key in column A, localized text in column B ---
$ws.Range("A6").Value = "Text/11/text"
$ws.Range("B6").Value = "Oh… It’s you… Have you heard what happened?"

$ws.Range("A7").Value = "Text/12/text"
$ws.Range("B7").Value = "I’m afraid… It came back…"

$ws.Range("A8").Value = "Text/13/text"
$ws.Range("B8").Value = "Yes… Yes… "

$ws.Range("A9").Value = "Text/14/text"
$ws.Range("B9").Value = "I can feel it… It’s getting close…"

# --- Highlight the "key" column for the localization table (green) ---
# RGB(0,169,51) foreground, RGB(0,128,0) background/pattern colour
$ws.Range("A2").Interior.Color = 3385600
$ws.Range("A2").Interior.PatternColor = 32768

$ws.Range("A3").Interior.Color = 3385600
$ws.Range("A3").Interior.PatternColor = 32768

$ws.Range("A4").Interior.Color = 3385600
$ws.Range("A4").Interior.PatternColor = 32768

$ws.Range("A5").Interior.Color = 3385600
$ws.Range("A5").Interior.PatternColor = 32768

$ws.Range("H2").Interior.Color = 3385600
$ws.Range("H2").Interior.PatternColor = 32768

# "strings" label cell - purple, RGB(128,0,128) fore/background
$ws.Range("H3").Interior.Color = 8388736
$ws.Range("H3").Interior.PatternColor = 8388736

# --- Column widths: new column A for keys, minor tweak to B (C/D untouched) ---
$ws.Columns.Item(1).ColumnWidth = 20.86
$ws.Columns.Item(2).ColumnWidth = 15.25

# --- Selection position, as recorded when the edit was made ---
$ws.Range("D11").Select()
